$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the question rows (A2:I26) by the "difficulté" column (I), ascending,
# keeping the header row (row 1) fixed. Sorting the full row range lets the
# row heights (used for the longer-answer rows) travel with their data.
$sortRange = $ws.Range("A2:I26")
$key1 = $ws.Range("I2:I26")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.SortMethod = 1
$ws.Sort.Apply()

# The "N°" column is a fixed running index, not part of the sorted data -
# restore it to 1..25 after the sort.
for ($i = 0; $i -lt 25; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $i + 1
}

# Move the active selection to E16, matching the recorded sheet view state.
$ws.Range("E16").Select()
